# Budget update - Resistor was added
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budget")

# Insert a new row at row 13, shifting the "Calcium hydroxyde" row (and everything
# below it) down by one. Formatting/formulas copy down from the row above, matching
# how the other line items in this table are built.
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the Resistor line item.
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "Resistor 3kOhms"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = "u"
$ws.Range("E13").Value = 8.41
$ws.Range("F13").Formula = "=+E13*C13"

$ws.Hyperlinks.Add(
    $ws.Range("G13"),
    "https://www.mouser.com/ProductDetail/Vishay-Dale/CMF073K0000JNEK?qs=w0S%252B%2FsJ%252B%2Fi0GZ3fUT3Vt%252BQ%3D%3D",
    [Type]::Missing,
    [Type]::Missing,
    "https://www.mouser.com/ProductDetail/Vishay-Dale/CMF073K0000JNEK?qs=w0S%252B%2FsJ%252B%2Fi0GZ3fUT3Vt%252BQ%3D%3D"
) | Out-Null
$ws.Range("G13").Style = "Hyperlink"

# The rows that shifted down (old rows 13-17, now 14-18) keep their own content but
# their sequential ID numbers (column A) need to bump by one.
$ws.Range("A14").Value = 11
$ws.Range("A15").Value = 12
$ws.Range("A16").Value = 13
$ws.Range("A17").Value = 14
$ws.Range("A18").Value = 15

# The power supply row now also shows its source URL as plain text (not a hyperlink).
$ws.Range("G12").Value = "https://www.amazon.com/%EF%BC%88Precision-00-01V%EF%BC%8C0-001A%EF%BC%894-Digital-Precision-Adjustable-Regulated/dp/B07M6JJS93"
